$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style of Q4 (year header) into R4, set the new year value
$ws.Range("R4").Value = 2021
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("R5").Value = 102.20441221981518
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R5").NumberFormat = "0.0"

$ws.Range("S9").Select()
